{"js": "// Remove the \"Word version of this document\" bullet item from the\n// \"Additional resources\" list (it pointed to a stale Word-doc download\n// link; the site now only links to PDF versions).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text.trim() === \"Word version of this document\") {\n    para.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Word version of this document\" bullet item from the\n# \"Additional resources\" list (the site now links to PDF versions instead\n# of a stand-alone Word-doc download).\n$d = $word.ActiveDocument\n\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd(\"`r\")\n    if ($t -eq \"Word version of this document\") {\n        $p.Range.Delete()\n    }\n}\n"}
